$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted at row 82, pushing the existing
# rows 82-128 down to 83-129 (dimension grows from A1:R128 to A1:R129).
$ws.Rows(82).Insert()

# Populate the newly inserted row with the latest weekly data point.
$ws.Range("A82").Value = 10
$ws.Range("B82").Value = "Vega Modelo de Temuco"
$ws.Range("C82").Value = "La Araucanía"
$ws.Range("D82").Value = 45236
$ws.Range("E82").Value = 9
$ws.Range("F82").Value = 100112022
$ws.Range("G82").Value = "Arveja Verde"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 150
$ws.Range("K82").Value = 21000
$ws.Range("L82").Value = 23000
$ws.Range("M82").Value = 21667
$ws.Range("N82").Value = "$/saco 25 kilos"
$ws.Range("O82").Value = "Región del Maule"
$ws.Range("P82").Value = 867
$ws.Range("Q82").Value = 25
$ws.Range("R82").Value = "Hortaliza"
